# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 17:35"

# Update province/community numeric data (columns B=Casos totales, C=Casos activos,
# D=Recuperados, E=Muertes). Only values that changed per the source diff are set.

# Row 4 - Madrid
$ws.Range("B4").Value = 68066
$ws.Range("D4").Value = 18639
$ws.Range("E4").Value = 8691

# Row 5 - Cataluña
$ws.Range("B5").Value = 57569
$ws.Range("D5").Value = 25783
$ws.Range("E5").Value = 5583

# Row 6 - Castilla y Leon
$ws.Range("B6").Value = 18591
$ws.Range("D6").Value = 7953
$ws.Range("E6").Value = 1922

# Row 7 - Castilla-La Mancha
$ws.Range("B7").Value = 16992
$ws.Range("D7").Value = 7655
$ws.Range("E7").Value = 2945

# Row 9 - Andalucia
$ws.Range("B9").Value = 12612
$ws.Range("D9").Value = 537
$ws.Range("E9").Value = 1404

# Row 16
$ws.Range("B16").Value = 5199
$ws.Range("D16").Value = 804

# Row 20
$ws.Range("B20").Value = 4044
$ws.Range("D20").Value = 577
$ws.Range("E20").Value = 360

# Row 32
$ws.Range("B32").Value = 2399
$ws.Range("D32").Value = 1026
$ws.Range("E32").Value = 310

# Row 33
$ws.Range("B33").Value = 2323
$ws.Range("D33").Value = 636

# Row 40
$ws.Range("B40").Value = 1587
$ws.Range("E40").Value = 148
